$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.950.55"
$ws.Range("E2").Value = "  -2.78%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.861.08"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.16"
$ws.Range("E5").Value = "  -1.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5028"
$ws.Range("E7").Value = "  -3.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3715"
$ws.Range("E8").Value = "  -1.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07118"
$ws.Range("E9").Value = "  -1.69%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8837"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.51"
$ws.Range("E11").Value = "  -2.89%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07566"
$ws.Range("E12").Value = "  -0.99%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.865.17"
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.285"
$ws.Range("E14").Value = "  -2.86%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "89.04"
$ws.Range("E15").Value = "  -3.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.29%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008366"
$ws.Range("E17").Value = "  -3.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.07"
$ws.Range("E18").Value = "  -2.82%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "26.993.33"
$ws.Range("E20").Value = "  -2.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.023"
$ws.Range("E21").Value = "  -2.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.112.73"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("E23").Value = "  -3.43%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.454"
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.848"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("E26").Value = "  -4.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.93"
$ws.Range("E27").Value = "  -2.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.097"
$ws.Range("E28").Value = "  -3.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "112.53"
$ws.Range("E29").Value = "  -1.96%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.637"
$ws.Range("E30").Value = "  -4.24%  "
$ws.Range("E31").Value = "  -3.11%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09024"
$ws.Range("E32").Value = "  +0.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05117"
$ws.Range("E33").Value = "  -3.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.045"
$ws.Range("E34").Value = "  -3.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.146"
$ws.Range("E35").Value = "  -7.62%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7216"
$ws.Range("E36").Value = "  -7.17%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02028"
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.029"
$ws.Range("E38").Value = "  -0.86%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.455"
$ws.Range("E39").Value = "  -6.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.072"
$ws.Range("E40").Value = "  -1.55%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5271"
$ws.Range("E41").Value = "  -4.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.502"
$ws.Range("E42").Value = "  -2.35%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "114.98"
$ws.Range("E43").Value = "  +1.51%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.228"
$ws.Range("E44").Value = "  -2.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1462"
$ws.Range("E45").Value = "  -3.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4591"
$ws.Range("E47").Value = "  -4.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.966"
$ws.Range("E48").Value = "  -3.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.560"
$ws.Range("E49").Value = "  -3.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.43"
$ws.Range("E50").Value = "  -1.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.76"
$ws.Range("E51").Value = "  -4.09%  "
